$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- userrole table (column D/E) ---
$ws.Range("D1").Value = "userrole"
$ws.Range("D3").Value = "id"
$ws.Range("D4").Value = "userid"
$ws.Range("D5").Value = "roleid"
$ws.Range("D6").Value = "raw_add_time"
$ws.Range("D7").Value = "raw_add_time"

# --- role table (column F/G) ---
$ws.Range("F1").Value = "role"
$ws.Range("F3").Value = "id"
$ws.Range("F4").Value = "name"
$ws.Range("F5").Value = "raw_add_time"
$ws.Range("F6").Value = "raw_add_time"

# --- rolepower table (column H/I) ---
$ws.Range("H1").Value = "rolepower"
$ws.Range("H3").Value = "id"
$ws.Range("H4").Value = "roleid"
$ws.Range("H5").Value = "powerid"
$ws.Range("H6").Value = "raw_add_time"
$ws.Range("H7").Value = "raw_add_time"

# --- power table (column J/K) ---
$ws.Range("J1").Value = "power"
$ws.Range("J3").Value = "id"
$ws.Range("J4").Value = "name"
$ws.Range("J5").Value = "raw_add_time"
$ws.Range("J6").Value = "raw_add_time"

# Match the column widths used by the existing A/B columns (13.875 chars)
$ws.Columns.Item(4).ColumnWidth = 13.14
$ws.Columns.Item(6).ColumnWidth = 13.14
$ws.Columns.Item(8).ColumnWidth = 13.14
$ws.Columns.Item(10).ColumnWidth = 13.14

# Restore the selection to match what was left after the edit
$ws.Range("J13").Select()
